$d = $word.ActiveDocument

# Update the date/day heading line
$d.Content.Find.Execute("2025-12-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-24 Wednesday", 2) | Out-Null

# Update the division problems in the single table, row by row (row 1, 5, 9, 13, 17
# contain the visible problems; the rows in between are blank spacer rows).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "17÷9="
$t.Cell(1,2).Range.Text = "96÷5="
$t.Cell(1,3).Range.Text = "31÷4="
$t.Cell(1,4).Range.Text = "56÷9="
$t.Cell(1,5).Range.Text = "63÷9="

$t.Cell(5,1).Range.Text = "26÷2="
$t.Cell(5,2).Range.Text = "38÷3="
$t.Cell(5,3).Range.Text = "38÷9="
$t.Cell(5,4).Range.Text = "24÷6="
$t.Cell(5,5).Range.Text = "45÷6="

$t.Cell(9,1).Range.Text = "37÷7="
$t.Cell(9,2).Range.Text = "19÷8="
$t.Cell(9,3).Range.Text = "60÷6="
$t.Cell(9,4).Range.Text = "92÷6="
$t.Cell(9,5).Range.Text = "33÷2="

$t.Cell(13,1).Range.Text = "54÷2="
$t.Cell(13,2).Range.Text = "77÷8="
$t.Cell(13,3).Range.Text = "62÷5="
$t.Cell(13,4).Range.Text = "76÷5="
$t.Cell(13,5).Range.Text = "53÷5="

$t.Cell(17,1).Range.Text = "96÷8="
$t.Cell(17,2).Range.Text = "65÷2="
$t.Cell(17,3).Range.Text = "49÷6="
$t.Cell(17,4).Range.Text = "20÷3="
$t.Cell(17,5).Range.Text = "87÷2="
